$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) The cells J1:K1 and N5 need to trade cell formats: J1/K1 (which
#    were formatted "center/center/wrap") must end up with the format
#    that N5 currently has ("vertical-center only"), and N5 must end
#    up with the format J1/K1 currently have. Stage J1's current
#    format in a scratch cell (Z1, well outside the used range) so it
#    survives the first overwrite, then swap, then clear the scratch
#    cell so it leaves no trace in the saved sheet.
# ------------------------------------------------------------------
$ws.Range("J1").Copy()
$ws.Range("Z1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("N5").Copy()
$ws.Range("J1:K1").PasteSpecial(-4122)

$ws.Range("Z1").Copy()
$ws.Range("N5").PasteSpecial(-4122)

$ws.Range("Z1").Clear()

# ------------------------------------------------------------------
# 2) Update the data values.
# ------------------------------------------------------------------
$ws.Range("L5").Value = 1.7
$ws.Range("N5").Value = 1.6

# ------------------------------------------------------------------
# 3) Remove the 2021 column (O) entirely - this drops O3/O4/O5 and
#    shrinks the sheet's used range/dimension and row spans from
#    column O back to N.
# ------------------------------------------------------------------
$ws.Columns("O:O").Delete()

# ------------------------------------------------------------------
# 4) Move the (inert) active-cell selection the way the author left it.
# ------------------------------------------------------------------
[void]$ws.Range("P6").Select()
